$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 442, shifting the existing rows 442:492
# down to 444:494.
$ws.Rows("442:443").Insert()

# Row 442 (new) - boilerplate columns shared with the rest of the block
$ws.Range("A442").Value = 7
$ws.Range("B442").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C442").Value = "Ñuble"
$ws.Range("D442").Value = 45142
$ws.Range("E442").Value = 16
$ws.Range("F442").Value = 100112003
$ws.Range("G442").Value = "Ajo"
$ws.Range("H442").Value = "Chino"
$ws.Range("I442").Value = "Primera"
$ws.Range("J442").Value = 30
$ws.Range("K442").Value = 21000
$ws.Range("L442").Value = 21000
$ws.Range("M442").Value = 21000
$ws.Range("N442").Value = "$/caja 10 kilos"
$ws.Range("O442").Value = "China"
$ws.Range("P442").Value = 2100
$ws.Range("Q442").Value = 10
$ws.Range("R442").Value = "Hortaliza"

# Row 443 (new)
$ws.Range("A443").Value = 7
$ws.Range("B443").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C443").Value = "Ñuble"
$ws.Range("D443").Value = 45142
$ws.Range("E443").Value = 16
$ws.Range("F443").Value = 100112003
$ws.Range("G443").Value = "Ajo"
$ws.Range("H443").Value = "Chino"
$ws.Range("I443").Value = "Primera"
$ws.Range("J443").Value = 30
$ws.Range("K443").Value = 23000
$ws.Range("L443").Value = 23000
$ws.Range("M443").Value = 23000
$ws.Range("N443").Value = "$/malla 10 kilos"
$ws.Range("O443").Value = "China"
$ws.Range("P443").Value = 2300
$ws.Range("Q443").Value = 10
$ws.Range("R443").Value = "Hortaliza"

# Make sure the date cells keep the date number format (style index 2 in the
# original file) that was already applied to column D by the Insert above.
$ws.Range("D442:D443").NumberFormat = "YYYY-MM-DD HH:MM:SS"
